$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11
$ws.Range("C2").Value = 1.006309338836623
$ws.Range("D2").Value = 0.3252074244777317

$ws.Range("C3").Value = 0.9588601731101325
$ws.Range("D3").Value = 0.3480492644806732

$ws.Range("C4").Value = 0.6707172461211615
$ws.Range("D4").Value = 0.5093814556948959

$ws.Range("C5").Value = 1.163850829661534
$ws.Range("D5").Value = 0.2569518331732663

$ws.Range("C6").Value = -0.05754694408945094
$ws.Range("D6").Value = 0.9546290025938191

$ws.Range("C7").Value = -0.5476907135780574
$ws.Range("D7").Value = 0.5894183856444006

$ws.Range("C8").Value = 0.09568540198998375
$ws.Range("D8").Value = 0.9246366887980015

$ws.Range("C9").Value = -0.3925403946015672
$ws.Range("D9").Value = 0.6984363573551868

$ws.Range("C10").Value = 0.1068164682965805
$ws.Range("D10").Value = 0.9159026672173893

$ws.Range("C11").Value = 0.5142032337409187
$ws.Range("D11").Value = 0.6122367218777476

$wb.Save()
